$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_test")

# The four "control" rows originally held unique per-occurrence labels
# ("Contrôle1", "Contrôle2", "Contrôle3", "Contrôle4") in column B.
# They should all read the same value as column C: "Contrôle".
$ws.Range("B12").Value = "Contrôle"
$ws.Range("B37").Value = "Contrôle"
$ws.Range("B38").Value = "Contrôle"
$ws.Range("B58").Value = "Contrôle"

# Widen column B to fit the new content (closest width this runtime's
# character-grid rounding can reach to the recorded 17.290714285714284).
$ws.Columns.Item(2).ColumnWidth = 16.5
